$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Update "Trends Status" sheet (sheet1)
# ---------------------------------------------------------------
$wsTrends = $wb.Worksheets.Item("Trends Status")
$wsTrends.Cells.Item(2,3).Value = 1      # C2
$wsTrends.Cells.Item(2,5).Value = 25     # E2
$wsTrends.Cells.Item(3,3).Value = 2      # C3
$wsTrends.Cells.Item(3,5).Value = 50     # E3
$wsTrends.Cells.Item(4,3).Value = 1      # C4
$wsTrends.Cells.Item(4,5).Value = 25     # E4
$wsTrends.Cells.Item(5,3).Value = 0      # C5
$wsTrends.Cells.Item(5,5).Value = 0      # E5
$wsTrends.Cells.Item(7,3).Value = 18     # C7
$wsTrends.Cells.Item(8,2).Value = 368    # B8
$wsTrends.Cells.Item(8,3).Value = 346    # C8

# ---------------------------------------------------------------
# 2. Update "Priority Status" sheet (sheet3)
# ---------------------------------------------------------------
$wsPriority = $wb.Worksheets.Item("Priority Status")
$wsPriority.Cells.Item(2,2).Value = 103  # B2
$wsPriority.Cells.Item(3,2).Value = 286  # B3
$wsPriority.Cells.Item(4,2).Value = 554  # B4

# ---------------------------------------------------------------
# 3. Update "Species qualification" sheet (sheet4)
# ---------------------------------------------------------------
$wsSpecies = $wb.Worksheets.Item("Species qualification")
$wsSpecies.Cells.Item(2,1).Value = "SoIB Assessment"  # A2
$wsSpecies.Cells.Item(2,2).Value = 368                # B2
$wsSpecies.Cells.Item(4,3).Value = 4                  # C4

# ---------------------------------------------------------------
# 4. Duplicate the "High Priority break-up" sheet so that the
#    original content is preserved on a new sheet named
#    "Major update - High Priority ", placed right after it, and
#    then overwrite the original sheet's data with the new
#    "Interannual update - High Pri" content, finally renaming it.
# ---------------------------------------------------------------
$wsHighPriority = $wb.Worksheets.Item("High Priority break-up")
$wsHighPriority.Copy([System.Reflection.Missing]::Value, $wsHighPriority)
$wsMajorUpdate = $wb.Worksheets.Item($wsHighPriority.Index + 1)
$wsMajorUpdate.Name = "Major update - High Priority "

# Now replace the contents of the original sheet with the new data
# and rename it.
$wsHighPriority.Cells.Clear()

$wsHighPriority.Cells.Item(1,1).Value = "Break-up"
$wsHighPriority.Cells.Item(1,2).Value = "High Species (no.)"
$wsHighPriority.Cells.Item(1,3).Value = "High Species (perc.)"
$wsHighPriority.Cells.Item(1,4).Value = "New High Species (no.)"
$wsHighPriority.Cells.Item(1,5).Value = "New High Species (perc.)"

$wsHighPriority.Cells.Item(2,1).Value = "Trend New"
$wsHighPriority.Cells.Item(2,2).Value = 73
$wsHighPriority.Cells.Item(2,3).Value = 70.90000000000001
$wsHighPriority.Cells.Item(2,4).Value = 73
$wsHighPriority.Cells.Item(2,5).Value = 78.5

$wsHighPriority.Cells.Item(3,1).Value = "IUCN"
$wsHighPriority.Cells.Item(3,2).Value = 30
$wsHighPriority.Cells.Item(3,3).Value = 29.1
$wsHighPriority.Cells.Item(3,4).Value = 20
$wsHighPriority.Cells.Item(3,5).Value = 21.5

$wsHighPriority.Name = "Interannual update - High Pri"
